# Update CVD (Cassville Voluntary-turnover Data?) figures across the per-location
# sheets so automate_finance.qmd picks up the refreshed CVD files dynamically.
$wb = $excel.ActiveWorkbook

# --- Piedras Negras Fasco Mexico: Professional Voluntary Turnover (Commit/Forecast), L4 -> 0
$ws = $wb.Worksheets.Item("Piedras Negras Fasco Mexico")
$ws.Range("L4").Value = 0

# --- Tipp City Ohio: Professional Voluntary Turnover (Commit/Forecast), L4 -> 0.0833
$ws = $wb.Worksheets.Item("Tipp City Ohio")
$ws.Range("L4").Value = 0.0833

# --- Faridabad India: Manufacturing Voluntary Turnover, E7/E8/E9 -> 0.0776
$ws = $wb.Worksheets.Item("Faridabad India")
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# --- Piedras Negras Jakel Mexico: Manufacturing Voluntary Turnover, E2/E3 -> 0.0776
$ws = $wb.Worksheets.Item("Piedras Negras Jakel Mexico")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776

# --- Fort Wayne Indiana: Professional Voluntary Turnover (Commit/Forecast), L4 -> 0.0345
#     Internal Fill Rate (Commit/Forecast), L7 cleared (blank)
$ws = $wb.Worksheets.Item("Fort Wayne Indiana")
$ws.Range("L4").Value = 0.0345
$ws.Range("L7").ClearContents()

# --- Grafton Wisconsin: Professional Voluntary Turnover (Commit/Forecast), L4 -> 0.04
$ws = $wb.Worksheets.Item("Grafton Wisconsin")
$ws.Range("L4").Value = 0.04

# --- Lavergne Tennessee: Internal Fill Rate (Commit/Forecast), L7 cleared (blank)
$ws = $wb.Worksheets.Item("Lavergne Tennessee")
$ws.Range("L7").ClearContents()

# --- Manila Philippines: Professional Voluntary Turnover (Commit/Forecast), L4 -> 0.1
$ws = $wb.Worksheets.Item("Manila Philippines")
$ws.Range("L4").Value = 0.1

# --- Milwaukee Pmc Hq Wisconsin: Internal Fill Rate, E5/E6/E7 -> 0.6
#     Commit/Forecast row (row 7): L7 -> 0, M7:W7 -> 0.6
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("E5").Value = 0.6
$ws.Range("E6").Value = 0.6
$ws.Range("E7").Value = 0.6
$ws.Range("L7").Value = 0
$ws.Range("M7:W7").Value = 0.6

# --- Mississauga Canada: Internal Fill Rate (Commit/Forecast), L5 cleared (blank)
$ws = $wb.Worksheets.Item("Mississauga Canada")
$ws.Range("L5").ClearContents()
